$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1470919.6
$ws.Range("J17").Value = 1470919.6
$ws.Range("L17").Value = 4412758.800000001
$ws.Range("N17").Value = -4413094.800000001
# Row 62
$ws.Range("H62").Value = 2756.8965
$ws.Range("I62").Value = 2451.842
$ws.Range("J62").Value = 3336.5
$ws.Range("K62").Value = 2451.842
$ws.Range("L62").Value = 3336.5
$ws.Range("M62").Value = -1827.842
$ws.Range("N62").Value = -4584.5
# Row 65
$ws.Range("H65").Value = 2756.8965
$ws.Range("I65").Value = 2451.842
$ws.Range("J65").Value = 3336.5
$ws.Range("K65").Value = 12259.21
$ws.Range("L65").Value = 16682.5
$ws.Range("M65").Value = -9139.210000000001
$ws.Range("N65").Value = -22922.5
# Row 70
$ws.Range("H70").Value = 2763.95
$ws.Range("I70").Value = 1087.5
$ws.Range("J70").Value = 3881.5833
$ws.Range("K70").Value = 3262.5
$ws.Range("L70").Value = 11644.7499
$ws.Range("M70").Value = -2992.5
$ws.Range("N70").Value = -12184.7499
# Row 73
$ws.Range("H73").Value = 2763.95
$ws.Range("I73").Value = 1087.5
$ws.Range("J73").Value = 3881.5833
$ws.Range("K73").Value = 3262.5
$ws.Range("L73").Value = 11644.7499
$ws.Range("M73").Value = -2326.5
$ws.Range("N73").Value = -13516.7499
# Row 98
$ws.Range("H98").Value = 1056
$ws.Range("I98").Value = 1056
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1056
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 442
$ws.Range("N98").ClearContents()
# Row 113
$ws.Range("H113").Value = 11735.786
$ws.Range("I113").Value = 2998.75
$ws.Range("J113").Value = 15230.6
$ws.Range("K113").Value = 2998.75
$ws.Range("L113").Value = 15230.6
$ws.Range("M113").Value = 255.25
$ws.Range("N113").Value = -21738.6
# Row 122
$ws.Range("H122").Value = 1056
$ws.Range("I122").Value = 1056
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3168
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -718
$ws.Range("N122").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1485.3103
$ws.Range("I61").Value = 1459.6666
$ws.Range("J61").Value = 1608.4
$ws.Range("K61").Value = 1459.6666
$ws.Range("L61").Value = 1608.4
$ws.Range("M61").Value = -1247.6666
$ws.Range("N61").Value = -2032.4
# Row 74
$ws.Range("H74").Value = 2009.6
$ws.Range("I74").Value = 2264.5
$ws.Range("J74").Value = 990
$ws.Range("K74").Value = 2264.5
$ws.Range("L74").Value = 990
$ws.Range("M74").Value = -1390.5
$ws.Range("N74").Value = -2738
# Row 77
$ws.Range("H77").Value = 2009.6
$ws.Range("I77").Value = 2264.5
$ws.Range("J77").Value = 990
$ws.Range("K77").Value = 11322.5
$ws.Range("L77").Value = 4950
$ws.Range("M77").Value = -6954.5
$ws.Range("N77").Value = -13686
# Row 125
$ws.Range("H125").Value = 78357.5
$ws.Range("J125").Value = 78357.5
$ws.Range("L125").Value = 78357.5
$ws.Range("N125").Value = -88197.5
# Row 136
$ws.Range("H136").Value = 1485.3103
$ws.Range("I136").Value = 1459.6666
$ws.Range("J136").Value = 1608.4
$ws.Range("K136").Value = 4378.9998
$ws.Range("L136").Value = 4825.200000000001
$ws.Range("M136").Value = -1828.9998
$ws.Range("N136").Value = -9925.200000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2533.516
$ws.Range("I58").Value = 1931.7
$ws.Range("K58").Value = 1931.7
$ws.Range("M58").Value = -1728.7
# Row 98
$ws.Range("H98").Value = 64890
$ws.Range("J98").Value = 64890
$ws.Range("L98").Value = 64890
$ws.Range("N98").Value = -69382
# Row 99
$ws.Range("H99").Value = 6248.773
$ws.Range("I99").Value = 1683.5
$ws.Range("J99").Value = 14238
$ws.Range("K99").Value = 1683.5
$ws.Range("L99").Value = 14238
$ws.Range("M99").Value = -185.5
$ws.Range("N99").Value = -17234
# Row 122
$ws.Range("H122").Value = 1667552.4
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 5000757
$ws.Range("K122").Value = 2850
$ws.Range("L122").Value = 15002271
$ws.Range("M122").Value = -400
$ws.Range("N122").Value = -15007171
# Row 126
$ws.Range("H126").Value = 6248.773
$ws.Range("I126").Value = 1683.5
$ws.Range("J126").Value = 14238
$ws.Range("K126").Value = 5050.5
$ws.Range("L126").Value = 42714
$ws.Range("M126").Value = -2580.5
$ws.Range("N126").Value = -47654
# Row 132
$ws.Range("H132").Value = 3419.8125
$ws.Range("I132").Value = 2724.889
$ws.Range("K132").Value = 8174.667
$ws.Range("M132").Value = -5644.667
# Row 136
$ws.Range("H136").Value = 2533.516
$ws.Range("I136").Value = 1931.7
$ws.Range("K136").Value = 5795.1
$ws.Range("M136").Value = -3245.1

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 4677.8125
$ws.Range("I109").Value = 5310
$ws.Range("J109").Value = 4390.4546
$ws.Range("K109").Value = 15930
$ws.Range("L109").Value = 13171.3638
$ws.Range("M109").Value = -14890
$ws.Range("N109").Value = -15251.3638
# Row 113
$ws.Range("H113").Value = 625.7436
$ws.Range("I113").Value = 497.84616
$ws.Range("K113").Value = 1493.53848
$ws.Range("M113").Value = 676.4615200000001
# Row 131
$ws.Range("H131").Value = 2864.2104
$ws.Range("J131").Value = 2954
$ws.Range("L131").Value = 8862
$ws.Range("N131").Value = -18942
# Row 132
$ws.Range("H132").Value = 1170.5714
$ws.Range("I132").Value = 911.7692
$ws.Range("J132").Value = 1591.125
$ws.Range("K132").Value = 8205.9228
$ws.Range("L132").Value = 14320.125
$ws.Range("M132").Value = -5675.9228
$ws.Range("N132").Value = -19380.125

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 30500
$ws.Range("J44").Value = 30500
$ws.Range("L44").Value = 30500
$ws.Range("N44").Value = -31692
# Row 46
$ws.Range("H46").Value = 6689.375
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 7073.5713
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 7073.5713
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -7385.5713
# Row 102
$ws.Range("H102").Value = 2040.6666
$ws.Range("I102").Value = 2064.6365
$ws.Range("J102").Value = 1777
$ws.Range("K102").Value = 2064.6365
$ws.Range("L102").Value = 1777
$ws.Range("M102").Value = -442.6365000000001
$ws.Range("N102").Value = -5021
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 126
$ws.Range("H126").Value = 2015.8276
$ws.Range("I126").Value = 1589.1177
$ws.Range("J126").Value = 2620.3333
$ws.Range("K126").Value = 4767.3531
$ws.Range("L126").Value = 7860.999899999999
$ws.Range("M126").Value = -2297.3531
$ws.Range("N126").Value = -12800.9999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 130976.25
$ws.Range("I40").Value = 500000
$ws.Range("J40").Value = 7968.3335
$ws.Range("K40").Value = 500000
$ws.Range("L40").Value = 7968.3335
$ws.Range("M40").Value = -499864
$ws.Range("N40").Value = -8240.333500000001
# Row 45
$ws.Range("H45").Value = 10000
$ws.Range("J45").Value = 10000
$ws.Range("L45").Value = 10000
$ws.Range("N45").Value = -10814
# Row 48
$ws.Range("H48").Value = 30000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 30000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 30000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -31322
# Row 122
$ws.Range("H122").Value = 4922.909
$ws.Range("I122").Value = 3657.7144
$ws.Range("J122").Value = 5513.3335
$ws.Range("K122").Value = 10973.1432
$ws.Range("L122").Value = 16540.0005
$ws.Range("M122").Value = -8523.143199999999
$ws.Range("N122").Value = -21440.0005

